$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row formatting for the two new rows (26, 27) from the closest
# existing rows that already carry the right visual style (border/wrap),
# then overwrite the values. ---
$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A26:E26").PasteSpecial(-4122) | Out-Null

$ws.Range("A24:E24").Copy() | Out-Null
$ws.Range("A27:E27").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# --- Row heights for the new rows ---
$ws.Rows.Item(26).RowHeight = 45
$ws.Rows.Item(27).RowHeight = 30

# --- Row heights changed on existing rows 19 and 20 ---
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 90

# --- New row 26 values (order A, C, B, D so shared-strings end up in the
# same order as the target workbook) ---
$ws.Range("A26").Value = "ENWIAM41"
$ws.Range("C26").Value = "Verify that the User is able to see 'Did you know? ...' Modal is displayed when user navigates from neon to ENW if Neon user has email same as existing steam acount"
$ws.Range("B26").Value = "OPQA-2172 || OPQA-1859"
$ws.Range("D26").Value = "N"

# --- New row 27 values ---
$ws.Range("A27").Value = "ENWIAM50"
$ws.Range("C27").Value = "Verify that Neon Landing page, displays Neon branding and marketing copy and also integration with Endnote"
$ws.Range("B27").Value = "OPQA-1707"
$ws.Range("D27").Value = "Y"

# --- Column B widened to fit new, longer content ---
$ws.Columns.Item(2).ColumnWidth = 23.5

# --- Selection moved to D26 ---
$ws.Range("D26").Select() | Out-Null
